# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.396.59'
$ws.Range("E2").Value = '  +0.22%  '

# Row 3
$ws.Range("D3").Value = '2.556.78'
$ws.Range("E3").Value = '  -2.45%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").Value = '''593.23'
$ws.Range("E5").Value = '  +0.75%  '

# Row 6
$ws.Range("D6").Value = '''173.57'
$ws.Range("E6").Value = '  +4.45%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").Value = '''0.531'
$ws.Range("E8").Value = '  -0.09%  '

# Row 9
$ws.Range("D9").Value = '2.556.14'
$ws.Range("E9").Value = '  -2.44%  '

# Row 10
$ws.Range("E10").Value = '  -0.48%  '

# Row 11
$ws.Range("E11").Value = '  +1.82%  '

# Row 12
$ws.Range("B12").Value = 'Toncoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D12").Value = '''5.17'
$ws.Range("E12").Value = '  -1.07%  '

# Row 13
$ws.Range("B13").Value = 'Cardano'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D13").Value = '''0.352'
$ws.Range("E13").Value = '  -4.53%  '

# Row 14
$ws.Range("D14").Value = '''27.15'
$ws.Range("E14").Value = '  -0.93%  '

# Row 15
$ws.Range("D15").Value = '3.017.55'
$ws.Range("E15").Value = '  -2.69%  '

# Row 16
$ws.Range("E16").Value = '  -0.80%  '

# Row 17
$ws.Range("D17").Value = '67.196.83'
$ws.Range("E17").Value = '  -0.23%  '

# Row 18
$ws.Range("D18").Value = '2.552.89'
$ws.Range("E18").Value = '  -2.74%  '

# Row 19
$ws.Range("D19").Value = '''8.07'
$ws.Range("E19").Value = '  +3.30%  '

# Row 20
$ws.Range("D20").Value = '''11.42'
$ws.Range("E20").Value = '  -2.95%  '

# Row 21
$ws.Range("D21").Value = '''356.86'
$ws.Range("E21").Value = '  +0.09%  '

# Row 22
$ws.Range("D22").Value = '''4.23'
$ws.Range("E22").Value = '  -1.60%  '

# Row 23
$ws.Range("D23").Value = '''4.69'
$ws.Range("E23").Value = '  +0.44%  '

# Row 24
$ws.Range("E24").Value = '  +5.29%  '

# Row 25
$ws.Range("E25").Value = '  -0.02%  '

# Row 26
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").Value = '''10.23'
$ws.Range("E26").Value = '  -3.57%  '

# Row 27
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").Value = '''70.22'
$ws.Range("E27").Value = '  +0.89%  '

# Row 28
$ws.Range("E28").Value = '  -2.76%  '

# Row 29
$ws.Range("E29").Value = '  +0.01%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0999'
$ws.Range("E30").Value = '  +0.02%  '

# Row 31
$ws.Range("D31").Value = '''539.48'
$ws.Range("E31").Value = '  -1.46%  '

# Row 32
$ws.Range("D32").Value = '''8.22'
$ws.Range("E32").Value = '  +1.56%  '

# Row 33
$ws.Range("E33").Value = '  +3.29%  '

# Row 34
$ws.Range("E34").Value = '  -0.27%  '

# Row 35
$ws.Range("E35").Value = '  -1.06%  '

# Row 36
$ws.Range("E36").Value = '  +0.01%  '

# Row 37
$ws.Range("E37").Value = '  +0.24%  '

# Row 38
$ws.Range("D38").Value = '''158.64'
$ws.Range("E38").Value = '  +0.55%  '

# Row 39
$ws.Range("D39").Value = '''18.81'
$ws.Range("E39").Value = '  -0.72%  '

# Row 40
$ws.Range("D40").Value = '''18.47'
$ws.Range("E40").Value = '  +1.15%  '

# Row 41
$ws.Range("E41").Value = '  -1.64%  '

# Row 42
$ws.Range("E42").Value = '  -0.06%  '

# Row 43
$ws.Range("D43").Value = '''5.20'
$ws.Range("E43").Value = '  +0.79%  '

# Row 44
$ws.Range("D44").Value = '''2.54'
$ws.Range("E44").Value = '  +4.53%  '

# Row 45
$ws.Range("E45").Value = '  -0.01%  '

# Row 46
$ws.Range("D46").Value = '''39.73'
$ws.Range("E46").Value = '  -1.22%  '

# Row 47
$ws.Range("D47").Value = '''151.20'
$ws.Range("E47").Value = '  -0.09%  '

# Row 48
$ws.Range("E48").Value = '  -2.25%  '

# Row 49
$ws.Range("E49").Value = '  -5.57%  '

# Row 50
$ws.Range("E50").Value = '  -1.43%  '

# Row 51
$ws.Range("D51").Value = '''1.72'
$ws.Range("E51").Value = '  +0.45%  '
